$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.771.07"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "1.976.53"
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'244.89"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").Value = "'0.630"
$ws.Range("E6").Value = "  +2.20%  "

$ws.Range("D7").Value = "'61.07"
$ws.Range("E7").Value = "  +3.74%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.380"
$ws.Range("E9").Value = "  +1.90%  "

$ws.Range("D10").Value = "'0.0794"
$ws.Range("E10").Value = "  -1.31%  "

$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("D12").Value = "'14.37"
$ws.Range("E12").Value = "  +4.94%  "

$ws.Range("D13").Value = "'0.845"
$ws.Range("E13").Value = "  +2.76%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'21.99"
$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.264.02"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").Value = "'5.42"
$ws.Range("E16").Value = "  +3.11%  "

$ws.Range("D17").Value = "1.965.53"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").Value = "36.646.75"
$ws.Range("E18").Value = "  +0.28%  "

$ws.Range("D19").Value = "'69.94"
$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("D20").Value = "0.0₃0859"
$ws.Range("E20").Value = "  +0.07%  "

$ws.Range("D21").Value = "'5.13"
$ws.Range("E21").Value = "  +1.24%  "

$ws.Range("D22").Value = "'230.00"
$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("D24").Value = "'2.47"
$ws.Range("E24").Value = "  +2.31%  "

$ws.Range("D25").Value = "'2.38"
$ws.Range("E25").Value = "  +1.97%  "

$ws.Range("D26").Value = "'0.145"
$ws.Range("E26").Value = "  +4.67%  "

$ws.Range("D27").Value = "'9.24"
$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("D28").Value = "'162.81"
$ws.Range("E28").Value = "  +1.64%  "

$ws.Range("D29").Value = "'19.45"
$ws.Range("E29").Value = "  +0.57%  "

$ws.Range("E30").Value = "  +21.10%  "

$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("D32").Value = "'4.83"
$ws.Range("E32").Value = "  +3.09%  "

$ws.Range("D33").Value = "'0.0620"
$ws.Range("E33").Value = "  +0.72%  "

$ws.Range("D34").Value = "'4.52"
$ws.Range("E34").Value = "  +6.24%  "

$ws.Range("E35").Value = "  +2.06%  "

$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("E37").Value = "  -1.65%  "

$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("D39").Value = "'5.44"
$ws.Range("E39").Value = "  -10.69%  "

$ws.Range("D40").Value = "'0.0974"
$ws.Range("E40").Value = "  -3.81%  "

$ws.Range("E41").Value = "  +0.63%  "

$ws.Range("E42").Value = "  +0.89%  "

$ws.Range("D43").Value = "'0.0211"
$ws.Range("E43").Value = "  -0.26%  "

$ws.Range("D44").Value = "'16.03"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").Value = "1.368.77"
$ws.Range("E45").Value = "  +0.77%  "

$ws.Range("D46").Value = "'89.49"
$ws.Range("E46").Value = "  +2.44%  "

$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("D48").Value = "'7.21"
$ws.Range("E48").Value = "  +1.24%  "

$ws.Range("D49").Value = "'2.83"
$ws.Range("E49").Value = "  -0.19%  "

$ws.Range("D50").Value = "'46.24"
$ws.Range("E50").Value = "  +6.58%  "

$ws.Range("D51").Value = "2.157.62"
$ws.Range("E51").Value = "  +0.77%  "
